$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'243.57"
$ws.Range("E2").Value = "'-0.86%"
$ws.Range("G2").Value = "'9"

# Row 3
$ws.Range("D3").Value = "'26.96"
$ws.Range("E3").Value = "'3.98%"
$ws.Range("G3").Value = "'9"

# Row 4
$ws.Range("D4").Value = "'5.157"
$ws.Range("E4").Value = "'0.81%"
$ws.Range("G4").Value = "'9"

# Row 5
$ws.Range("D5").Value = "'0.05618"
$ws.Range("E5").Value = "'0.46%"
$ws.Range("G5").Value = "'9"

# Row 6
$ws.Range("D6").Value = "'6.491"
$ws.Range("G6").Value = "'9"

# Row 7
$ws.Range("D7").Value = "'0.8167"
$ws.Range("E7").Value = "'-0.01%"
$ws.Range("G7").Value = "'9"

# Row 8
$ws.Range("D8").Value = "'0.8308"
$ws.Range("E8").Value = "'-1.98%"
$ws.Range("G8").Value = "'9"

# Row 9
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.009939"
$ws.Range("E9").Value = "'1,556.35%"
$ws.Range("G9").Value = "'9"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1326"
$ws.Range("E10").Value = "'-1.17%"
$ws.Range("G10").Value = "'9"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.06927"
$ws.Range("E11").Value = "'-0.29%"
$ws.Range("G11").Value = "'9"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02890"
$ws.Range("E12").Value = "'1.32%"
$ws.Range("G12").Value = "'9"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09375"
$ws.Range("E13").Value = "'-0.32%"
$ws.Range("G13").Value = "'9"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001524"
$ws.Range("E14").Value = "'-0.06%"
$ws.Range("G14").Value = "'9"

# Row 15
$ws.Range("D15").Value = "'0.006157"
$ws.Range("E15").Value = "'1.12%"
$ws.Range("G15").Value = "'9"

# Row 16
$ws.Range("D16").Value = "'3.610"
$ws.Range("E16").Value = "'1.72%"
$ws.Range("G16").Value = "'9"

# Row 17
$ws.Range("D17").Value = "'3.023"
$ws.Range("E17").Value = "'-0.24%"
$ws.Range("G17").Value = "'9"

# Row 18
$ws.Range("D18").Value = "'2.304"
$ws.Range("E18").Value = "'8.76%"
$ws.Range("G18").Value = "'9"

# Row 19
$ws.Range("E19").Value = "'-0.72%"
$ws.Range("G19").Value = "'9"

# Row 20
$ws.Range("D20").Value = "'0.03088"
$ws.Range("E20").Value = "'-4.52%"
$ws.Range("G20").Value = "'9"

# Row 21
$ws.Range("E21").Value = "'-2.21%"
$ws.Range("G21").Value = "'9"

# Row 22
$ws.Range("E22").Value = "'0.01%"
$ws.Range("G22").Value = "'9"

# Row 23
$ws.Range("D23").Value = "'0.04584"
$ws.Range("E23").Value = "'-2.57%"
$ws.Range("G23").Value = "'9"

# Row 24
$ws.Range("E24").Value = "'-2.48%"
$ws.Range("G24").Value = "'9"

# Row 25
$ws.Range("D25").Value = "'0.001226"
$ws.Range("E25").Value = "'-1.84%"
$ws.Range("G25").Value = "'9"

# Row 26
$ws.Range("D26").Value = "'0.004485"
$ws.Range("E26").Value = "'-2.68%"
$ws.Range("G26").Value = "'9"

# Row 27
$ws.Range("G27").Value = "'9"

# Row 28
$ws.Range("E28").Value = "'0.67%"
$ws.Range("G28").Value = "'9"

# Row 29
$ws.Range("G29").Value = "'9"

# Row 30
$ws.Range("G30").Value = "'9"

# Row 31
$ws.Range("G31").Value = "'9"

# Row 32
$ws.Range("G32").Value = "'9"

# Row 33
$ws.Range("G33").Value = "'9"

# Row 34
$ws.Range("G34").Value = "'9"

# Row 35
$ws.Range("G35").Value = "'9"

# Row 36
$ws.Range("G36").Value = "'9"

# Row 37
$ws.Range("G37").Value = "'9"

# Row 38
$ws.Range("G38").Value = "'9"

# Row 39
$ws.Range("G39").Value = "'9"

# Row 40
$ws.Range("E40").Value = "'-0.47%"
$ws.Range("G40").Value = "'9"

# Row 41
$ws.Range("D41").Value = "'0.006080"
$ws.Range("E41").Value = "'-0.82%"
$ws.Range("G41").Value = "'9"

# Row 42
$ws.Range("D42").Value = "'0.1051"
$ws.Range("E42").Value = "'-0.18%"
$ws.Range("G42").Value = "'9"

# Row 43
$ws.Range("D43").Value = "'0.002573"
$ws.Range("E43").Value = "'2.90%"
$ws.Range("G43").Value = "'9"

# Row 44
$ws.Range("D44").Value = "'0.008316"
$ws.Range("E44").Value = "'5.35%"
$ws.Range("G44").Value = "'9"

# Row 45
$ws.Range("D45").Value = "'0.00005298"
$ws.Range("E45").Value = "'-0.37%"
$ws.Range("G45").Value = "'9"

# Row 46
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("G46").Value = "'9"

# Row 47
$ws.Range("E47").Value = "'-18.36%"
$ws.Range("G47").Value = "'9"

# Row 48
$ws.Range("D48").Value = "'0.002639"
$ws.Range("E48").Value = "'23.99%"
$ws.Range("G48").Value = "'9"

# Row 49
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("G49").Value = "'9"

# Row 50
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("G50").Value = "'9"

# Row 51
$ws.Range("G51").Value = "'9"
